# Update Data by bot, scripted by HH
# Row 6 (003029 / 吉大正元) is refreshed from the 2019 annual-report figures
# to the 2020 Q3 (three-quarter) report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date
$ws.Range("H6").Value = "2020-09-30 00:00:00"

# Numeric financial metrics
$ws.Range("I6").Value = 0.21
$ws.Range("J6").Value = 0.21
$ws.Range("K6").Value = 283341662.91
$ws.Range("L6").Value = 27951921.58
$ws.Range("M6").Value = 3.99
$ws.Range("N6").Value = -2.7527442875
$ws.Range("O6").Value = 3.7762566403
$ws.Range("P6").Value = 5.281553560015
$ws.Range("Q6").Value = -0.322052988248
$ws.Range("R6").Value = 62.3293607217

# S6 / T6 no longer carry a value for this report period - blank them out
# but keep them as (empty) text cells, matching the sibling placeholder
# cells (U6, V6, X6, ...) already present on this row.
$ws.Range("S6").Value = "'"
$ws.Range("S6").Style = "Normal"
$ws.Range("T6").Value = "'"
$ws.Range("T6").Style = "Normal"

# Flags / labels describing which report this row represents.
# These columns are stored as text even though they look numeric, so force
# a text number format before writing, then drop back to the default style
# so no stray formatting is left behind on the cell.
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "1"
$ws.Range("AB6").Style = "Normal"

$ws.Range("AC6").Value = "2020Q3"
$ws.Range("AD6").Value = "2020年 三季报"

$ws.Range("AE6").NumberFormat = "@"
$ws.Range("AE6").Value = "2020"
$ws.Range("AE6").Style = "Normal"

$ws.Range("AF6").Value = "三季报"

# Timestamp this refresh was produced
$ws.Range("AG6").Value = "2020-12-07 07:57:14"
